# pairedGSEA metadata: final analysis + de-duplicate the "LMNA deficiency"
# comparison titles into distinct per-comparison labels.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "LMNA deficiency - control"
$ws.Range("G3").Value = "LMNA deficiency - donor"

# Match the author's final selection/view state (active cell G4).
$ws.Range("G4").Select()
